$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "OFF Pre-op"
$ws.Range("B1").Value = "OFF Post-op"
$ws.Range("C1").Value = "DCS-P2"

$ws.Range("A2").Value = 0.13941390229079501
$ws.Range("B2").Value = 0.20930321573712601
$ws.Range("C2").Value = 0.196557193564451

$ws.Range("A3").Value = 0.085462824872357901
$ws.Range("B3").Value = 0.146737267056747
$ws.Range("C3").Value = 0.14762375977915301

$ws.Range("A4").Value = 0.072806262996169893
$ws.Range("B4").Value = 0.12784250595379601
$ws.Range("C4").Value = 0.179867071228175

$ws.Range("A5").Value = 0.087909183305049196
$ws.Range("B5").Value = 0.120380478682331
$ws.Range("C5").Value = 0.13952919185812801

$ws.Range("A6").Value = 0.077039934990282402
$ws.Range("B6").Value = 0.16514477906925101
$ws.Range("C6").Value = 0.13103633343806301

$ws.Range("A7").Value = 0.123901976163654
$ws.Range("B7").Value = 0.17174324658212001
$ws.Range("C7").Value = 0.11155736866136

$ws.Range("B8").Value = 0.073145008425080996
$ws.Range("B9").Value = 0.098551642935804998
$ws.Range("B10").Value = 0.10279716489537701
$ws.Range("B11").Value = 0.083531123792990297
$ws.Range("B12").Value = 0.110125720772694
$ws.Range("B13").Value = 0.19008619760360801

$ws.Range("A2:C7").Font.Name = "Arial"
$ws.Range("A2:C7").Font.Size = 10
$ws.Range("A8:B13").Font.Name = "Arial"
$ws.Range("A8:B13").Font.Size = 10

[void]$ws.Range("F10").Select()
